$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 0.4269757926895136
$ws.Range("J3").Value = 0.5451654426702852
$ws.Range("K3").Value = 0.5531585209632127
$ws.Range("L3").Value = 2.81890317602702
